# Insert one new data row into the "Feria Lagunitas de Puerto Montt - Cebolla" sheet.
# This pushes the existing rows 311..379 down to 312..380 and fills the
# freshly-inserted row 311 with a new weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 311:379 down by one row, creating a blank row 311.
$ws.Rows("311:311").Insert()

# Populate the new row 311 with the new observation.
$ws.Range("A311").Value = 4
$ws.Range("B311").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C311").Value = "Los Lagos"
$ws.Range("D311").Value = 44543
$ws.Range("E311").Value = 10
$ws.Range("F311").Value = 100112004
$ws.Range("G311").Value = "Cebolla"
$ws.Range("H311").Value = "Sin especificar"
$ws.Range("I311").Value = "Primera"
$ws.Range("J311").Value = 250
$ws.Range("K311").Value = 8000
$ws.Range("L311").Value = 8000
$ws.Range("M311").Value = 8000
$ws.Range("N311").Value = "`$/malla 18 kilos"
$ws.Range("O311").Value = "Perú"
$ws.Range("P311").Value = 444
$ws.Range("Q311").Value = 18
$ws.Range("R311").Value = "Hortaliza"
